$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set a cell's text value while preventing Excel from auto-converting
# numeric-looking strings (like price values) into floating point numbers.
function Set-TextValue {
    param($Cell, $Text)
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
}

Set-TextValue "D2" "244.85"
Set-TextValue "D4" "5.327"
Set-TextValue "D5" "0.05787"
Set-TextValue "D6" "6.481"
Set-TextValue "D7" "3.330"
Set-TextValue "D8" "0.8107"
Set-TextValue "D9" "0.8878"
Set-TextValue "D10" "0.1395"
Set-TextValue "D11" "0.07352"
Set-TextValue "D13" "0.03054"
Set-TextValue "D14" "0.09328"
Set-TextValue "D15" "3.875"
Set-TextValue "D16" "0.001538"
Set-TextValue "D17" "0.04711"
Set-TextValue "D18" "0.0006026"
Set-TextValue "D19" "0.006049"
Set-TextValue "D20" "0.001295"
Set-TextValue "D21" "0.00008806"
Set-TextValue "D24" "0.3179"
Set-TextValue "D27" "0.004616"
Set-TextValue "D28" "0.0002351"
Set-TextValue "D41" "0.006364"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.004103"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1054"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue "D44" "0.007615"
Set-TextValue "D45" "0.00005471"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.5504"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.0002001"
